$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column I: "APTO?" formula ---
# Row 6 gets its own (non-shared) formula
$ws.Range("I6").Formula = '=IF(OR(AND(C6="Madrid",H6=1,G6<22),AND(H6=2,E6="M",G6>18),AND(F6="ADMINISTRATIVO",E6="F",C6="Barcelona",D6="COMPLETO"),AND(E6="F",D6="VESPERTINO",C6="Madrid")),"APTO","NO APTO")'

# Rows 7-17 share one formula
$ws.Range("I7:I17").Formula = '=IF(OR(AND(C7="Madrid",H7=1,G7<22),AND(H7=2,E7="M",G7>18),AND(F7="ADMINISTRATIVO",E7="F",C7="Barcelona",D7="COMPLETO"),AND(E7="F",D7="VESPERTINO",C7="Madrid")),"APTO","NO APTO")'

# --- Column J: name if criteria met, else "SIN NOMBRE" ---
$ws.Range("J6").Formula = '=IF(OR(AND(D6="COMPLETO",F6="ADMINISTRATIVO",C6="MADRID"),AND(G6=19),OR(G6=20),AND(D6="COMPLETO"),OR(D6="MATUTINO")),B6,"SIN NOMBRE")'

$ws.Range("J7:J17").Formula = '=IF(OR(AND(D7="COMPLETO",F7="ADMINISTRATIVO",C7="MADRID"),AND(G7=19),OR(G7=20),AND(D7="COMPLETO"),OR(D7="MATUTINO")),B7,"SIN NOMBRE")'

# --- Q12: a blank-line note, wrapped, which grows row 12's height ---
$ws.Range("Q12").Value = "
"
$ws.Range("Q12").WrapText = $true
$ws.Rows.Item(12).RowHeight = 30

# --- Column J width (best-fit for the new names) ---
$ws.Columns.Item(10).ColumnWidth = 13.6

# --- Selection left where the author's cursor ended up ---
[void]$ws.Range("J21").Select()
